$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1380.8
$ws.Range("I43").Value = 977.25
$ws.Range("K43").Value = 977.25
$ws.Range("M43").Value = -908.25

$ws.Range("H54").Value = 11557.333
$ws.Range("I54").Value = 14538
$ws.Range("K54").Value = 14538
$ws.Range("M54").Value = -14052

$ws.Range("H137").Value = 3391.7778
$ws.Range("I137").Value = 2170.8572
$ws.Range("J137").Value = 3943.1614
$ws.Range("K137").Value = 6512.571599999999
$ws.Range("L137").Value = 11829.4842
$ws.Range("M137").Value = -3962.571599999999
$ws.Range("N137").Value = -16929.4842

$ws.Range("H138").Value = 931934.7
$ws.Range("I138").Value = 2428.1365
$ws.Range("J138").Value = 1201002.4
$ws.Range("K138").Value = 7284.4095
$ws.Range("L138").Value = 3603007.2
$ws.Range("M138").Value = -2144.4095
$ws.Range("N138").Value = -3613287.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13568.593
$ws.Range("I32").Value = 10754.652
$ws.Range("K32").Value = 10754.652
$ws.Range("M32").Value = -10467.652

$ws.Range("H58").Value = 80000
$ws.Range("J58").Value = 80000
$ws.Range("L58").Value = 80000
$ws.Range("N58").Value = -80860

$ws.Range("H61").Value = 10190.552
$ws.Range("I61").Value = 5635.9443
$ws.Range("J61").Value = 17643.545
$ws.Range("K61").Value = 5635.9443
$ws.Range("L61").Value = 17643.545
$ws.Range("M61").Value = -5423.9443
$ws.Range("N61").Value = -18067.545

$ws.Range("H122").Value = 31253250
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550

$ws.Range("H132").Value = 3790.6936
$ws.Range("I132").Value = 1199.6945
$ws.Range("K132").Value = 3599.0835
$ws.Range("M132").Value = -1069.0835

$ws.Range("H136").Value = 10190.552
$ws.Range("I136").Value = 5635.9443
$ws.Range("J136").Value = 17643.545
$ws.Range("K136").Value = 16907.8329
$ws.Range("L136").Value = 52930.63499999999
$ws.Range("M136").Value = -14357.8329
$ws.Range("N136").Value = -58030.63499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20626

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 13217.777
$ws.Range("J4").Value = 13217.777
$ws.Range("L4").Value = 13217.777
$ws.Range("N4").Value = -13441.777

$ws.Range("H31").Value = 3905.0952
$ws.Range("I31").Value = 1001.0833
$ws.Range("J31").Value = 7777.1113
$ws.Range("K31").Value = 1001.0833
$ws.Range("L31").Value = 7777.1113
$ws.Range("M31").Value = -706.0833
$ws.Range("N31").Value = -8367.1113

$ws.Range("H34").Value = 3905.0952
$ws.Range("I34").Value = 1001.0833
$ws.Range("J34").Value = 7777.1113
$ws.Range("K34").Value = 1001.0833
$ws.Range("L34").Value = 7777.1113
$ws.Range("M34").Value = -799.0833
$ws.Range("N34").Value = -8181.1113

$ws.Range("H58").Value = 2167477.5
$ws.Range("I58").Value = 3248974
$ws.Range("J58").Value = 4484.2856
$ws.Range("K58").Value = 3248974
$ws.Range("L58").Value = 4484.2856
$ws.Range("M58").Value = -3248771
$ws.Range("N58").Value = -4890.2856

$ws.Range("H102").Value = 46000
$ws.Range("J102").Value = 46000
$ws.Range("L102").Value = 46000
$ws.Range("N102").Value = -50868

$ws.Range("H122").Value = 18431.77
$ws.Range("I122").Value = 10633.333
$ws.Range("J122").Value = 25116.143
$ws.Range("K122").Value = 31899.999
$ws.Range("L122").Value = 75348.429
$ws.Range("M122").Value = -29449.999
$ws.Range("N122").Value = -80248.429

$ws.Range("H136").Value = 2167477.5
$ws.Range("I136").Value = 3248974
$ws.Range("J136").Value = 4484.2856
$ws.Range("K136").Value = 9746922
$ws.Range("L136").Value = 13452.8568
$ws.Range("M136").Value = -9744372
$ws.Range("N136").Value = -18552.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1890.909
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6224

$ws.Range("H113").Value = 772.0465
$ws.Range("I113").Value = 778.5
$ws.Range("J113").Value = 760
$ws.Range("K113").Value = 2335.5
$ws.Range("L113").Value = 2280
$ws.Range("M113").Value = -165.5
$ws.Range("N113").Value = -6620

$ws.Range("H131").Value = 588.14435
$ws.Range("I131").Value = 307.15686
$ws.Range("J131").Value = 899.6739
$ws.Range("K131").Value = 921.4705799999999
$ws.Range("L131").Value = 2699.0217
$ws.Range("M131").Value = 4118.52942
$ws.Range("N131").Value = -12779.0217

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12300
$ws.Range("J5").Value = 12300
$ws.Range("L5").Value = 12300
$ws.Range("N5").Value = -12524

$ws.Range("H97").Value = 2075.111
$ws.Range("I97").Value = 2270.7693
$ws.Range("J97").Value = 1566.4
$ws.Range("K97").Value = 2270.7693
$ws.Range("L97").Value = 1566.4
$ws.Range("M97").Value = -1774.7693
$ws.Range("N97").Value = -2558.4

$ws.Range("H132").Value = 42107.07
$ws.Range("I132").Value = 64636.562
$ws.Range("J132").Value = 12067.75
$ws.Range("K132").Value = 193909.686
$ws.Range("L132").Value = 36203.25
$ws.Range("M132").Value = -191379.686
$ws.Range("N132").Value = -41263.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 385000
$ws.Range("J2").Value = 94444.44500000001
$ws.Range("L2").Value = 94444.44500000001
$ws.Range("N2").Value = -94668.44500000001

$ws.Range("H22").Value = 1693.2
$ws.Range("I22").Value = 1574.75
$ws.Range("J22").Value = 1772.1666
$ws.Range("K22").Value = 1574.75
$ws.Range("L22").Value = 1772.1666
$ws.Range("M22").Value = -1279.75
$ws.Range("N22").Value = -2362.1666

$ws.Range("H27").Value = 1693.2
$ws.Range("I27").Value = 1574.75
$ws.Range("J27").Value = 1772.1666
$ws.Range("K27").Value = 1574.75
$ws.Range("L27").Value = 1772.1666
$ws.Range("M27").Value = -1467.75
$ws.Range("N27").Value = -1986.1666

$ws.Range("H60").Value = 20061
$ws.Range("J60").Value = 20061
$ws.Range("L60").Value = 20061
$ws.Range("N60").Value = -21079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2889.111
$ws.Range("I2").Value = 2002
$ws.Range("K2").Value = 2002
$ws.Range("M2").Value = -1890

$ws.Range("H136").Value = 8630
$ws.Range("I136").Value = 9755.556
$ws.Range("J136").Value = 8268.214
$ws.Range("K136").Value = 29266.668
$ws.Range("L136").Value = 24804.642
$ws.Range("M136").Value = -26716.668
$ws.Range("N136").Value = -29904.642
